# Apply updated dSF (column F) values per repull of data / mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -7
    "F3"  = -6
    "F5"  = 3
    "F10" = -3
    "F13" = 1
    "F18" = 1
    "F19" = 1
    "F20" = -1
    "F25" = -4
    "F26" = -1
    "F29" = 1
    "F32" = -4
    "F33" = 3
    "F38" = -4
    "F40" = 2
    "F41" = -3
    "F44" = 0
    "F49" = -2
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

$wb.Save()
